$wb = $excel.ActiveWorkbook

# --- DecisionTree sheet: append rows 31-35 (copies of row 22), fix last row's I/J ---
$ws3 = $wb.Worksheets.Item("DecisionTree")
$srcRow3 = $ws3.Range("A22:AK22").Value()
for ($r = 31; $r -le 35; $r++) {
    $ws3.Range("A$r" + ":AK$r").Value = $srcRow3
}
$ws3.Range("I35").Value = $ws3.Range("W22").Value()
$ws3.Range("J35").Value = $ws3.Range("Z22").Value()

# --- NaiveBayes sheet: append rows 10-14 (copies of row 9), fix last row's D/E ---
$ws4 = $wb.Worksheets.Item("NaiveBayes")
$srcRow4 = $ws4.Range("A9:AF9").Value()
for ($r = 10; $r -le 14; $r++) {
    $ws4.Range("A$r" + ":AF$r").Value = $srcRow4
}
$ws4.Range("D14").Value = $ws4.Range("R9").Value()
$ws4.Range("E14").Value = $ws4.Range("U9").Value()

# --- LogisticRegression sheet: append rows 10-14 (copies of row 9), fix last row's E/F ---
$ws5 = $wb.Worksheets.Item("LogisticRegression")
$srcRow5 = $ws5.Range("A9:AG9").Value()
for ($r = 10; $r -le 14; $r++) {
    $ws5.Range("A$r" + ":AG$r").Value = $srcRow5
}
$ws5.Range("E14").Value = $ws5.Range("S9").Value()
$ws5.Range("F14").Value = $ws5.Range("V9").Value()
